$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 27: TFX pipeline post -> React Native deployment track post
$ws.Range("D27").Value = "React Native 앱의 배포 트랙 관리하기"
$ws.Range("E27").Value = "https://tech.scatterlab.co.kr/react-native-deployment-track/"

# Row 42: python error handling post -> 위경도 거리 차이 post
$ws.Range("D42").Value = "위경도 도분초, 십진법 값별 거리 차이"
$ws.Range("E42").Value = "https://kjk92.tistory.com/96"

# Row 50: openCV super-resolution post -> 정보기하학과 머신러닝 post
$ws.Range("D50").Value = "정보기하학과 머신러닝"
$ws.Range("E50").Value = "http://incredible.egloos.com/7558260"

# Row 51: [독후감] economic book review -> [python] zero-filled list post
$ws.Range("D51").Value = "[python] 0으로 채워진 1차원, 2차원 리스트(배열) 만들기"
$ws.Range("E51").Value = "https://bskyvision.com/entry/python-0%EC%9C%BC%EB%A1%9C-%EC%B1%84%EC%9B%8C%EC%A7%84-1%EC%B0%A8%EC%9B%90-2%EC%B0%A8%EC%9B%90-%EB%A6%AC%EC%8A%A4%ED%8A%B8%EB%B0%B0%EC%97%B4-%EB%A7%8C%EB%93%A4%EA%B8%B0"
